$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the two records (by their ID in column A) that this edit removes
# entirely from the sample.
$rm232 = $ws.Columns.Item(1).Find("RM 232")
$sc92  = $ws.Columns.Item(1).Find("SC 92")

# Delete the lower row first so the row number already captured for the
# upper row stays valid.
if ($rm232.Row -lt $sc92.Row) {
    $ws.Rows.Item($sc92.Row).Delete()
    $ws.Rows.Item($rm232.Row).Delete()
} else {
    $ws.Rows.Item($rm232.Row).Delete()
    $ws.Rows.Item($sc92.Row).Delete()
}

# With those two rows gone, re-locate the remaining records whose "C" column
# (column D) missingness changes: "SC 5" now has its true value revealed,
# while "SC 101" becomes the newly-missing record.
$sc5 = $ws.Columns.Item(1).Find("SC 5")
$ws.Cells.Item($sc5.Row, 4).Value = -13.8

$sc101 = $ws.Columns.Item(1).Find("SC 101")
$ws.Cells.Item($sc101.Row, 4).Value = ""
